$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos ")

# Header renames (G1 then F1) -- also updates Table1 column names automatically
$ws.Range("G1").Value = " "
$ws.Range("F1").Value = "Columna1"

# New summary block title row
$ws.Range("F14").Value = "REQUERIMIENTO 2"
$ws.Range("G14").Value = "REQUERIMIENTO 4"

# Labels down column E
$ws.Range("E15").Value = "Small:"
$ws.Range("E16").Value = "10pct:"
$ws.Range("E17").Value = "50pct:"
$ws.Range("E18").Value = "large:"

# Values down column F
$ws.Range("F15").Value = "4.96ms"
$ws.Range("F16").Value = "10.3ms"
$ws.Range("F17").Value = "30.44ms"
$ws.Range("F18").Value = "40.52ms"

# Values down column G
$ws.Range("G15").Value = "6.00ms"
$ws.Range("G16").Value = "15.647ms"
$ws.Range("G17").Value = "43.678ms"
$ws.Range("G18").Value = "90,68ms"

# Row 14 height grows to fit the wrapped header text
$ws.Rows.Item(14).RowHeight = 28.5

# Column widths widen to fit new content
$ws.Columns.Item(6).ColumnWidth = 19.75
$ws.Columns.Item(7).ColumnWidth = 19.5

# Selection / scroll position as last left by the author
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F26").Select()
